{"js": "// Add \"no border\" tblBorders (top/left/bottom/right/insideH -- note: no\n// insideV) to the \"Personal Information\" table (the table that holds the\n// \"Cultural language\" / \"Laz language\" row) so the grid lines Word would\n// otherwise draw from the TableGrid style are suppressed.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Locate the \"Personal Information\" table by looking for the row whose\n// label cell reads \"Cultural language\" (robust even if table ordering\n// ever changes).\nfor (let i = 0; i < tables.items.length; i++) {\n  tables.items[i].load(\"values\");\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  const values = tables.items[i].values;\n  const hit = values.some((row) => row.some((cell) => cell.indexOf(\"Cultural language\") !== -1));\n  if (hit) {\n    target = tables.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  // Fallback: third table in document order (0-based index 2) matches the\n  // \"Personal Information\" table in this r\u00e9sum\u00e9 layout.\n  target = tables.items[2];\n}\n\nconst borderLocations = [\"Top\", \"Left\", \"Bottom\", \"Right\", \"InsideHorizontal\"];\nfor (const location of borderLocations) {\n  const border = target.getBorder(location);\n  border.width = 0;\n  border.type = \"None\";\n}\n\nawait context.sync();\n", "ps1": "# Add \"no border\" tblBorders (top/left/bottom/right/insideH -- note: no\n# insideV) to the \"Personal Information\" table (the table that holds the\n# \"Cultural language\" / \"Laz language\" row) so the grid lines Word would\n# otherwise draw from the TableGrid style are suppressed.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Personal Information\" table by looking for the row whose\n# label cell reads \"Cultural language\" (robust even if table ordering\n# ever changes).\n$target = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    if ($t.Range.Text -like \"*Cultural language*\") {\n        $target = $t\n        break\n    }\n}\n\nif ($target -eq $null) {\n    # Fallback: third table in document order matches the \"Personal\n    # Information\" table in this resume layout.\n    $target = $d.Tables.Item(3)\n}\n\n$borders = $target.Borders\n\n# wdBorderTop=-1, wdBorderLeft=-2, wdBorderBottom=-3, wdBorderRight=-4,\n# wdBorderHorizontal=-5 (insideH). wdBorderVertical (-6 / insideV) is\n# intentionally left untouched.\n$locations = @(-1, -2, -3, -4, -5)\nforeach ($location in $locations) {\n    $border = $borders.Item($location)\n    $border.LineWidth = 0\n    $border.LineStyle = 0\n}\n"}
